# Fruta / hortaliza, semanal
# Insert 3 new daily price rows for Pera (Packham's Triumph / Winter Nelis)
# right before the existing row 477, pushing the rest of the table down by
# three rows (old A1:T570 -> new A1:T573).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new records: shift rows 477-570 down to 480-573.
$ws.Rows("477:479").Insert()

# Common columns (Mercado/Region/etc.) are identical across this sheet.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100104
$producto   = "Frutos de pepita"
$categoriaId = 100104005
$categoria   = "Pera"

function Set-Row($r, $fecha, $variedad, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $origen, $pkg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $pkg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 477 44505 "Packham's Triumph" "Especial" 650 15000 16000 15538 "$/caja 18 kilos granel" "Provincia de Curicó" 863 18
Set-Row 478 44505 "Packham's Triumph" "Primera"  470 13000 14000 13468 "$/caja 18 kilos granel" "Provincia de Curicó" 748 18
Set-Row 479 44505 "Winter Nelis"      "Primera"  300 14000 14000 14000 "$/caja 18 kilos granel" "Región de O'Higgins" 778 18
